$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "94÷9=10, 4" "43÷2=21, 1"
Replace-Text "42÷8=5, 2" "89÷6=14, 5"
Replace-Text "70÷3=23, 1" "12÷6=2, 0"
Replace-Text "74÷9=8, 2" "41÷2=20, 1"
Replace-Text "24÷7=3, 3" "57÷7=8, 1"
Replace-Text "69÷9=7, 6" "99÷8=12, 3"
Replace-Text "72÷6=12, 0" "36÷2=18, 0"
Replace-Text "52÷5=10, 2" "75÷2=37, 1"
Replace-Text "90÷6=15, 0" "97÷6=16, 1"
Replace-Text "32÷8=4, 0" "56÷4=14, 0"
Replace-Text "34÷2=17, 0" "77÷9=8, 5"
Replace-Text "81÷3=27, 0" "63÷6=10, 3"
Replace-Text "74÷2=37, 0" "34÷9=3, 7"
Replace-Text "84÷8=10, 4" "99÷9=11, 0"
Replace-Text "60÷7=8, 4" "19÷9=2, 1"
Replace-Text "65÷6=10, 5" "83÷9=9, 2"
Replace-Text "61÷7=8, 5" "15÷9=1, 6"
Replace-Text "36÷9=4, 0" "70÷9=7, 7"
Replace-Text "16÷6=2, 4" "68÷9=7, 5"
Replace-Text "52÷3=17, 1" "10÷7=1, 3"
Replace-Text "91÷7=13, 0" "29÷7=4, 1"
Replace-Text "39÷3=13, 0" "10÷2=5, 0"
Replace-Text "63÷5=12, 3" "28÷2=14, 0"
Replace-Text "28÷4=7, 0" "37÷6=6, 1"
Replace-Text "19÷2=9, 1" "24÷8=3, 0"

Write-Host "Done replacing all values"
